# Updates crypto price (D) and volume-change (E) columns to refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.649.73"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.344.50"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.23"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.16"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -8.52%  "
$ws.Range("D9").Value = "2.342.82"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.49"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "2.768.61"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "60.566.74"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "2.344.58"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.71"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.26"
$ws.Range("E26").Value = "  +8.21%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "498.61"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").Value = "0.0₃0856"
$ws.Range("E30").Value = "  -8.52%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.374"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.21"
$ws.Range("E38").Value = "  -4.62%  "
$ws.Range("E39").Value = "  +5.56%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.57"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  -6.65%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.35"
$ws.Range("E51").Value = "  -1.58%  "
